$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Terrible's - 101)
$ws.Range("C2").Value = 16.6
$ws.Range("E2").Value = 97.90000000000001
$ws.Range("F2").Value = 49

# Row 4 (Terrible's - 118)
$ws.Range("B4").Value = 252
$ws.Range("C4").Value = 99.59999999999999
$ws.Range("D4").Value = 22.3
$ws.Range("E4").Value = 373.9
$ws.Range("F4").Value = 187

# Row 7 (Terrible's - 129)
$ws.Range("B7").Value = 93.8
$ws.Range("C7").Value = 39.6
$ws.Range("D7").Value = 39.8
$ws.Range("E7").Value = 173.2
$ws.Range("F7").Value = 87

# Row 8 (Terrible's - 132)
$ws.Range("B8").Value = 92.40000000000001
$ws.Range("E8").Value = 154.8
$ws.Range("F8").Value = 78

# Row 19 (Terrible's - 156)
$ws.Range("B19").Value = 35.6
$ws.Range("C19").Value = 11.8
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 62.4
$ws.Range("F19").Value = 31

# Row 21 (Terrible's - 162)
$ws.Range("B21").Value = 297.4
$ws.Range("C21").Value = 87.40000000000001
$ws.Range("D21").Value = 35.3
$ws.Range("E21").Value = 421.1
$ws.Range("F21").Value = 210

# Row 23 (Terrible's - 165)
$ws.Range("B23").Value = 185
$ws.Range("C23").Value = 61.6
$ws.Range("D23").Value = 64.5
$ws.Range("E23").Value = 311.1
$ws.Range("F23").Value = 156

# Row 26 (Terrible's - 170)
$ws.Range("E26").Value = 163.3
$ws.Range("F26").Value = 82
